$d = $word.ActiveDocument

# The paragraph "I soon find out that, ... sees me." used to be split across
# three separate <w:r> runs ("I soon " / "find" / " out that ... sees me.").
# Re-typing the whole sentence over itself collapses it back down to a
# single run containing the full, unbroken text.
$apostrophe = [char]0x2019
$sentence = "I soon find out that, despite not being late, I still made her wait a bit. However, instead of being put-out like a certain pink-haired childhood friend might" + $apostrophe + "ve been, she instead waves when she sees me."

$found = $d.Content.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)
if (-not $found) {
    throw "Could not locate the 'I soon find out...' sentence to merge its runs."
}
